# "Updated values for NZE" - update NZE scenario input values on the
# "Data" sheet. Column EW holds the 2050 end-point for each scenario row;
# columns DT:EV are shared formulas ( previous_cell + ($EW-$DS)/30 ) that
# interpolate from the DS (start, column 123) value out to the EW (end)
# value, so they recompute automatically once EW is written.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# New NZE end-values (column EW) for the affected scenario rows.
$ws.Range("EW2").Value  = 0.8
$ws.Range("EW5").Value  = 1.3
$ws.Range("EW6").Value  = 0.8
$ws.Range("EW10").Value = 0.5
$ws.Range("EW15").Value = 1.1
$ws.Range("EW16").Value = 1
$ws.Range("EW17").Value = 1

# Move the active selection on the frozen/split pane to EW7, matching
# where the author was last working in the sheet.
$ws.Range("EW7").Select() | Out-Null
